$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 290
$ws.Range("C3").Value = 179743
$ws.Range("C4").Value = 169699
$ws.Range("C8").Value = 65.02
